$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 600
$ws.Range("F7").Value = 14976
$ws.Range("F8").Value = 405
$ws.Range("F10").Value = 677
$ws.Range("F11").Value = 15229
$ws.Range("F13").Value = 8747
$ws.Range("F22").Value = 513
$ws.Range("F23").Value = 22
$ws.Range("F25").Value = 51
$ws.Range("F33").Value = 31
$ws.Range("F37").Value = 427
$ws.Range("F39").Value = 5373

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 600
$ws.Range("F7").Value = 14976
$ws.Range("F8").Value = 405
$ws.Range("F10").Value = 677
$ws.Range("F11").Value = 15229
$ws.Range("F13").Value = 8747
$ws.Range("F23").Value = 513
$ws.Range("F24").Value = 22
$ws.Range("F26").Value = 51
$ws.Range("F36").Value = 31
$ws.Range("F40").Value = 427
$ws.Range("F42").Value = 5373

$wb.Save()
